$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E7").Value = 11.82109999999999
$ws.Range("A8").Value = -21.15330000000001
$ws.Range("A10").Value = -20.50559999999997
$ws.Range("A12").Value = -22.43200000000003
$ws.Range("C13").Value = -12.89849999999999
$ws.Range("A18").Value = -22.36120000000003
$ws.Range("E20").Value = 12.1727
